# Update session 39 documents.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- Title: "Assignment by Monday, Dec 4" -> "Assignment by Wednesday, Dec 6"
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Assignment by Wednesday, Dec 6"

# --- Content placeholder: reorder / update the bullet paragraphs.
# Paragraph 1 ("Assignment:") and the two trailing blank paragraphs stay untouched.
# Before: 2) Continue..., 3) Cloud Computing..., 4) Signup...
# After:  2) Verify...,    3) Signup...,          4) Continue..., 5) Review Cloud...
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$tr.Paragraphs(2,1).Text = "Verify the you can access Final Exam Study Guide… let’s do that now"
$tr.Paragraphs(3,1).Text = "Signup for a Project 5 review timeslot for Friday… if you are not already presenting in class"
$tr.Paragraphs(4,1).Text = "Continue to focus on Project 5… which is due one week from today"
$inserted = $tr.Paragraphs(4,1).InsertAfter("`rReview Cloud Computing and MapReduce slides")
